$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in G1 header: remove duplicate "类"
$ws.Range("G1").Value = "鞋类居民消费价格指数(上年=100)"

# Copy formatting of the existing "year" cell (A6) down into the two new
# year rows (A7, A8) before filling in the values, so the bold/border/
# centered style used for the year column is preserved.
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)

# Add new row 7 (2021年)
$ws.Range("A7").Value = "2021年"
$ws.Range("D7").Value = 100.4
$ws.Range("F7").Value = 100.3
$ws.Range("G7").Value = 99.8

# Add new row 8 (2022年)
$ws.Range("A8").Value = "2022年"
$ws.Range("F8").Value = 100.5
